# modal summary correction v5
#
# The invoice extraction columns L:T (issueDate, AFM_issuer, AA, series,
# totalNetValue, totalVatAmount, totalValue, character, type) were
# duplicated/extraneous metadata columns and are removed.
# The first data row (old row 2, MARK 400011172559639) was a stray/duplicate
# summary entry and is removed, leaving only the second invoice
# (MARK 400011184530011) as the single data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra metadata columns L through T (used to be A1:T3).
$ws.Range("L1:T3").EntireColumn.Delete() | Out-Null

# Remove the first data row; the second invoice row shifts up to row 2.
$ws.Range("A2:K2").EntireRow.Delete() | Out-Null

Write-Host "Final UsedRange:" $ws.UsedRange.Address()
